# Fix: "Time in solve_pbtk is days not hours!"
#
# The underlying TK/QSPR model was refit after correcting the time units
# bug, which changed the raw regression statistics stored on the
# 'Cmax-stats' and 'AUC-stats' sheets. 'Table4' only contains
# INDEX/MATCH formulas pulling from those two sheets, so updating the
# raw values here and recalculating is sufficient to refresh Table4 too.

$wb = $excel.ActiveWorkbook

$cmax = $wb.Worksheets.Item("Cmax-stats")
$auc  = $wb.Worksheets.Item("AUC-stats")

# --- Cmax-stats (columns: A=QSPR, B=R2, C=RMSLE, D=RPE) ---
$cmaxValues = @{
    2 = @(0.567, 0.836, -0.0229)
    3 = @(0.621, 0.782, 0.15)
    4 = @(0.596, 0.808, 0.0169)
    5 = @(0.57, 0.833, -0.105)
    6 = @(0.578, 0.826, 0.111)
    8 = @(0.464, 0.93, 0.531)
    9 = @(0.604, 0.864, 0.156)
}

foreach ($row in $cmaxValues.Keys) {
    $vals = $cmaxValues[$row]
    $cmax.Cells.Item($row, 2).Value = $vals[0]
    $cmax.Cells.Item($row, 3).Value = $vals[1]
    $cmax.Cells.Item($row, 4).Value = $vals[2]
}

# --- AUC-stats (columns: A=QSPR, B=R2, C=RMSLE, D=RPE, E=RPE.low, F=RPE.high) ---
$aucValues = @{
    2 = @(0.495, 1.11, 3.44, 6.72, -0.464)
    3 = @(0.618, 0.966, 5.88, 10.5, -0.437)
    4 = @(0.584, 1.01, 3.32, 5.0, -0.0382)
    5 = @(0.396, 1.22, 2.29, 5.99, -0.796)
    6 = @(0.56, 1.04, 4.67, 7.93, 0.385)
    7 = @(0.961, 0.309, 0.315, 0.411, 0.0618)
    8 = @(0.102, 1.48, 5.1, 12.1, -0.881)
    9 = @(0.622, 1.08, 4.65, 10.5, -0.118)
}

foreach ($row in $aucValues.Keys) {
    $vals = $aucValues[$row]
    $auc.Cells.Item($row, 2).Value = $vals[0]
    $auc.Cells.Item($row, 3).Value = $vals[1]
    $auc.Cells.Item($row, 4).Value = $vals[2]
    $auc.Cells.Item($row, 5).Value = $vals[3]
    $auc.Cells.Item($row, 6).Value = $vals[4]
}

# Match the reviewer's selection state left on each stats sheet after
# refreshing the numbers (cosmetic, but mirrors the authored change).
$cmax.Activate()
$cmax.Range("B2:D9").Select()

$auc.Activate()
$auc.Range("B2:F9").Select()

$excel.Calculate()
